# Change Friday killSp (column L) from 30 to 20 for rows 27-31
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("InvasionTable")

$ws.Range("L27:L31").Value = 20
